# Insert two new observation rows above the current first data row (row 2),
# pushing all existing rows down by two, then populate the two new rows
# with the new species-observation data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 2 (existing rows 2-8 become rows 4-10).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# --- New row 2 ---
$ws.Range("A2").Value = 111908383
$ws.Range("B2").Value = 90670
$ws.Range("C2").Value = "Ovaliderad"
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 4365
$ws.Range("F2").Value = "Smalfotad taggsvamp"
$ws.Range("G2").Value = "Hydnellum gracilipes"
$ws.Range("H2").Value = "(P.Karst) P.Karst"
$ws.Range("P2").Value = "Prästtjärnen (Prästtjärnen), Dlr"
$ws.Range("Q2").Value = 518003.61510633
$ws.Range("R2").Value = 6789983.610409672
$ws.Range("S2").Value = 25
$ws.Range("T2").Value = "Dalarna"
$ws.Range("U2").Value = "Rättvik"
$ws.Range("V2").Value = "Dalarna"
$ws.Range("W2").Value = "Ore"
$ws.Range("Y2").Value = "'2023-09-05"
$ws.Range("Z2").Value = "'13:25"
$ws.Range("AA2").Value = "'2023-09-05"
$ws.Range("AB2").Value = "'13:25"
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false
$ws.Range("AW2").Value = "Andreas Öster"
$ws.Range("AX2").Value = "Andreas Öster"

# --- New row 3 ---
$ws.Range("A3").Value = 111908386
$ws.Range("B3").Value = 88489
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 1962
$ws.Range("F3").Value = "Vaddporing"
$ws.Range("G3").Value = "Anomoporia kamtschatica"
$ws.Range("H3").Value = "(Parmasto) Bondartseva"
$ws.Range("P3").Value = "Prästtjärnen (Prästtjärnen), Dlr"
$ws.Range("Q3").Value = 518003.61510633
$ws.Range("R3").Value = 6789983.610409672
$ws.Range("S3").Value = 25
$ws.Range("T3").Value = "Dalarna"
$ws.Range("U3").Value = "Rättvik"
$ws.Range("V3").Value = "Dalarna"
$ws.Range("W3").Value = "Ore"
$ws.Range("Y3").Value = "'2023-09-05"
$ws.Range("Z3").Value = "'13:25"
$ws.Range("AA3").Value = "'2023-09-05"
$ws.Range("AB3").Value = "'13:25"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = "Andreas Öster"
$ws.Range("AX3").Value = "Andreas Öster"
